$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously-empty Wireframe (column C) cells for existing Login rows ---
$ws.Range("C22").Value = "WireFrame_Client_001"
$ws.Range("C23").Value = "WireFrame_Login_005"
$ws.Range("C25").Value = "WireFrame_Login_002"

# --- Insert a new traceability row for the Login feature, right after row 25 ---
$ws.Range("A26").EntireRow.Insert()

# Copy formatting from the row above (25) onto the new row 26, then fix up the
# Wireframe cell (C26) to use the wrap-text style used elsewhere in the sheet (row 39, column C).
$ws.Range("A25:C25").Copy($ws.Range("A26:C26"))
$ws.Range("C39").Copy($ws.Range("C26"))

# Populate the new row's content.
$ws.Range("A26").Value = "SRS_Login_005"
$ws.Range("B26").Value = "TC_Login_002`nTC_Login_005`nTC_Login_006`nTC_Login_007`nTC_Login_0012`nTC_Login_0013`n"
$ws.Range("C26").Value = "WireFrame_Login_003`nWireFrame_Login_004`nWireFrame_Login_005"

# Row heights: the new row is taller (wrapped multi-line content); rows 27/28 (old 26/27)
# also end up with adjusted heights after the insert.
$ws.Rows("26:26").RowHeight = 42
$ws.Rows("27:27").RowHeight = 21
$ws.Rows("28:28").RowHeight = 31.5

# --- Restore the view state (selection / scroll position) seen in the edited workbook ---
$ws.Application.Goto($ws.Range("A20"))
$ws.Range("F26").Select()
